$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 4.33
$ws.Range("I2").Value = 4.5
$ws.Range("N2").Value = 19
$ws.Range("X2").Value = 11
$ws.Range("AH2").Value = 19
$ws.Range("AL2").Value = 29
$ws.Range("AQ2").Value = 23
$ws.Range("BC2").Value = 101
